$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1515.5
$ws.Range("I18").Value = 1515.5
$ws.Range("K18").Value = 1515.5
$ws.Range("M18").Value = -1231.5

$ws.Range("H34").Value = 8000
$ws.Range("I34").Value = 8000
$ws.Range("K34").Value = 8000
$ws.Range("M34").Value = -7797

$ws.Range("H36").Value = 8000
$ws.Range("I36").Value = 8000
$ws.Range("K36").Value = 8000
$ws.Range("M36").Value = -7285

$ws.Range("H38").Value = 423.77777
$ws.Range("I38").Value = 351.75
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 1055.25
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -683.25
$ws.Range("N38").Value = -3744

$ws.Range("H58").Value = 95
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H64").Value = 3199.5
$ws.Range("J64").Value = 3199.5
$ws.Range("L64").Value = 3199.5
$ws.Range("N64").Value = -3695.5

$ws.Range("H67").Value = 3199.5
$ws.Range("J67").Value = 3199.5
$ws.Range("L67").Value = 3199.5
$ws.Range("N67").Value = -4915.5

$ws.Range("H87").Value = 79999
$ws.Range("J87").Value = 79999
$ws.Range("L87").Value = 79999
$ws.Range("N87").Value = -82495

$ws.Range("H90").Value = 79999
$ws.Range("J90").Value = 79999
$ws.Range("L90").Value = 239997
$ws.Range("N90").Value = -252477

$ws.Range("H132").Value = 1868.125
$ws.Range("I132").Value = 1868.125
$ws.Range("K132").Value = 5604.375
$ws.Range("M132").Value = -3074.375

$ws.Range("H135").Value = 933.3226
$ws.Range("I135").Value = 593.2174
$ws.Range("J135").Value = 1911.125
$ws.Range("K135").Value = 5338.9566
$ws.Range("L135").Value = 17200.125
$ws.Range("M135").Value = -2803.9566
$ws.Range("N135").Value = -22270.125

$ws.Range("H138").Value = 2958.697
$ws.Range("I138").Value = 2091.6667
$ws.Range("J138").Value = 3999.1333
$ws.Range("K138").Value = 6275.000100000001
$ws.Range("L138").Value = 11997.3999
$ws.Range("M138").Value = -1135.000100000001
$ws.Range("N138").Value = -22277.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3464.4517
$ws.Range("I32").Value = 2635.0876
$ws.Range("K32").Value = 2635.0876
$ws.Range("M32").Value = -2348.0876

$ws.Range("H97").Value = 653.1111
$ws.Range("I97").Value = 633.86664
$ws.Range("J97").Value = 749.3333
$ws.Range("K97").Value = 633.86664
$ws.Range("L97").Value = 749.3333
$ws.Range("M97").Value = -137.86664
$ws.Range("N97").Value = -1741.3333

$ws.Range("H122").Value = 1441142.6
$ws.Range("I122").Value = 1680833
$ws.Range("K122").Value = 5042499
$ws.Range("M122").Value = -5040049

$ws.Range("H132").Value = 2458.4443
$ws.Range("I132").Value = 2446.7144
$ws.Range("K132").Value = 7340.1432
$ws.Range("M132").Value = -4810.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3856.7144
$ws.Range("I94").Value = 3499.25
$ws.Range("K94").Value = 3499.25
$ws.Range("M94").Value = -3048.25

$ws.Range("H105").Value = 3872.353
$ws.Range("I105").Value = 4473.857
$ws.Range("K105").Value = 4473.857
$ws.Range("M105").Value = -2726.857

$ws.Range("H134").Value = 3468.2144
$ws.Range("I134").Value = 1936.4286
$ws.Range("K134").Value = 5809.2858
$ws.Range("M134").Value = -3274.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 204.52942
$ws.Range("I7").Value = 147.3077
$ws.Range("K7").Value = 147.3077
$ws.Range("M7").Value = -34.30770000000001

$ws.Range("H22").Value = 363.33334
$ws.Range("I22").Value = 295
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 295
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 55
$ws.Range("N22").Value = -1200

$ws.Range("H107").Value = 585.75
$ws.Range("I107").Value = 338.4
$ws.Range("K107").Value = 338.4
$ws.Range("M107").Value = 1581.6

$ws.Range("H134").Value = 2690.9614
$ws.Range("I134").Value = 2527.7222
$ws.Range("J134").Value = 3058.25
$ws.Range("K134").Value = 7583.1666
$ws.Range("L134").Value = 9174.75
$ws.Range("M134").Value = -5048.1666
$ws.Range("N134").Value = -14244.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 99250
$ws.Range("J37").Value = 99250
$ws.Range("L37").Value = 297750
$ws.Range("N37").Value = -297974

$ws.Range("H132").Value = 3820.6924
$ws.Range("I132").Value = 4018.7778
$ws.Range("J132").Value = 3375
$ws.Range("K132").Value = 36169.00019999999
$ws.Range("L132").Value = 30375
$ws.Range("M132").Value = -33639.00019999999
$ws.Range("N132").Value = -35435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2316.2
$ws.Range("I31").Value = 395.25
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 395.25
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -103.25
$ws.Range("N31").Value = -10584

$ws.Range("H35").Value = 6333333.5
$ws.Range("I35").Value = 7000000
$ws.Range("K35").Value = 7000000
$ws.Range("M35").Value = -6999702

$ws.Range("H37").Value = 2316.2
$ws.Range("I37").Value = 395.25
$ws.Range("J37").Value = 10000
$ws.Range("K37").Value = 395.25
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = -118.25
$ws.Range("N37").Value = -10554

$ws.Range("H46").Value = 51972
$ws.Range("J46").Value = 34999.5
$ws.Range("L46").Value = 34999.5
$ws.Range("N46").Value = -35311.5

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 5768
$ws.Range("J80").Value = 6664
$ws.Range("L80").Value = 6664
$ws.Range("N80").Value = -8660

$ws.Range("H83").Value = 5768
$ws.Range("J83").Value = 6664
$ws.Range("L83").Value = 33320
$ws.Range("N83").Value = -43304

$ws.Range("H102").Value = 100000
$ws.Range("I102").Value = 100000
$ws.Range("K102").Value = 100000
$ws.Range("M102").Value = -98378

$ws.Range("H113").Value = 100000
$ws.Range("J113").Value = 100000
$ws.Range("L113").Value = 100000
$ws.Range("N113").Value = -104340

$ws.Range("H132").Value = 4199.4
$ws.Range("I132").Value = 3999
$ws.Range("K132").Value = 11997
$ws.Range("M132").Value = -9467

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1428.1428
$ws.Range("I7").Value = 1499.75
$ws.Range("K7").Value = 1499.75
$ws.Range("M7").Value = -1387.75

$ws.Range("H22").Value = 5250
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 7500
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 7500
$ws.Range("M22").Value = -2705
$ws.Range("N22").Value = -8090

$ws.Range("H27").Value = 5250
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 7500
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 7500
$ws.Range("M27").Value = -2893
$ws.Range("N27").Value = -7714

$ws.Range("H40").Value = 3921.75
$ws.Range("I40").Value = 2848.5
$ws.Range("J40").Value = 4995
$ws.Range("K40").Value = 2848.5
$ws.Range("L40").Value = 4995
$ws.Range("M40").Value = -2712.5
$ws.Range("N40").Value = -5267

$ws.Range("H55").Value = 478.08334
$ws.Range("I55").Value = 164.875
$ws.Range("J55").Value = 1104.5
$ws.Range("K55").Value = 164.875
$ws.Range("L55").Value = 1104.5
$ws.Range("M55").Value = 8.125
$ws.Range("N55").Value = -1450.5

$ws.Range("H68").Value = 3280.4
$ws.Range("I68").Value = 2199
$ws.Range("J68").Value = 4001.3333
$ws.Range("K68").Value = 2199
$ws.Range("L68").Value = 4001.3333
$ws.Range("M68").Value = -1450
$ws.Range("N68").Value = -5499.3333

$ws.Range("H71").Value = 3280.4
$ws.Range("I71").Value = 2199
$ws.Range("J71").Value = 4001.3333
$ws.Range("K71").Value = 10995
$ws.Range("L71").Value = 20006.6665
$ws.Range("M71").Value = -7251
$ws.Range("N71").Value = -27494.6665

$ws.Range("H93").Value = 1106.4667
$ws.Range("I93").Value = 810.5
$ws.Range("J93").Value = 1698.4
$ws.Range("K93").Value = 810.5
$ws.Range("L93").Value = 1698.4
$ws.Range("M93").Value = 437.5
$ws.Range("N93").Value = -4194.4

$ws.Range("H126").Value = 1428.1428
$ws.Range("I126").Value = 1499.75
$ws.Range("K126").Value = 4499.25
$ws.Range("M126").Value = -2029.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2215.6428
$ws.Range("I126").Value = 1963.8889
$ws.Range("K126").Value = 5891.6667
$ws.Range("M126").Value = -3421.6667
